# Update "想去人数" (column F) values for several events on both the
# "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Map event name (column C) -> new value for column F
$updates = @{
    "南宁·THOsp·幻想朱槿绘翠5"           = 113
    "南宁·原x穹x崩only"                 = 250
    "南宁·第五人格Only1.0"               = 123
    "南宁·AP动漫游戏嘉年华"              = 1611
    "南宁·布谷鸟动漫展4th"               = 1458
    "南宁·恋与深空only"                 = 263
    "南宁·小蜜蜂动漫嘉年华2.0"           = 56
    "南宁·AB动漫游戏嘉年华"              = 402
    "横州·第二届海棠动漫游戏嘉年华"       = 108
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastRow = $ws.UsedRange.Rows.Count

    for ($row = 2; $row -le $lastRow; $row++) {
        $name = $ws.Cells.Item($row, 3).Value2   # column C = 名称

        if (-not [string]::IsNullOrEmpty($name)) {
            if ($updates.ContainsKey($name)) {
                $ws.Cells.Item($row, 6).Value = $updates[$name]   # column F = 想去人数
            }
        }
    }
}
